$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row data: row, Coin(B), Link(C), Price(D), Volume1h(E)
$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "51.691.71", "  +4.08%  ")
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "2.774.22", "  +5.18%  ")
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.07%  ")
    ,@(5, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "116.45", "  +2.97%  ")
    ,@(6, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "333.59", "  +2.90%  ")
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.539", "  +1.88%  ")
    ,@(8, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "0.999", "  -0.14%  ")
    ,@(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.576", "  +5.62%  ")
    ,@(10, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "41.96", "  +5.25%  ")
    ,@(11, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.0865", "  +6.35%  ")
    ,@(12, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "20.35", "  +2.50%  ")
    ,@(13, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.130", "  +2.31%  ")
    ,@(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "7.67", "  +4.85%  ")
    ,@(15, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "3.203.54", "  +5.03%  ")
    ,@(16, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "2.774.90", "  +5.08%  ")
    ,@(17, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.892", "  +3.63%  ")
    ,@(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "51.652.85", "  +4.15%  ")
    ,@(19, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "3.31", "  +11.28%  ")
    ,@(20, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "13.57", "  +5.18%  ")
    ,@(21, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "6.87", "  +2.52%  ")
    ,@(22, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0₃0978", "  +3.14%  ")
    ,@(23, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "278.34", "  +3.02%  ")
    ,@(24, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "69.84", "  +1.21%  ")
    ,@(25, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "2.68", "  +5.44%  ")
    ,@(26, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "26.90", "  +2.07%  ")
    ,@(27, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.00", "  +0.07%  ")
    ,@(28, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "10.19", "  -1.40%  ")
    ,@(29, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "2.23", "  +0.22%  ")
    ,@(30, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.142", "  +1.71%  ")
    ,@(31, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "35.19", "  +0.31%  ")
    ,@(32, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "50.23", "  +1.23%  ")
    ,@(33, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "5.59", "  +1.87%  ")
    ,@(34, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.0821", "  +0.80%  ")
    ,@(35, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "1.00", "  -0.11%  ")
    ,@(36, "Celestia", "https://coinranking.com/coin/YQcD0lBl7+celestia-tia", "19.08", "  +0.39%  ")
    ,@(37, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "5.04", "  +2.07%  ")
    ,@(38, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "2.10", "  +2.38%  ")
    ,@(39, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "3.25", "  +4.03%  ")
    ,@(40, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0354", "  +9.04%  ")
    ,@(41, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "127.25", "  +0.32%  ")
    ,@(42, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "23.12", "  +2.55%  ")
    ,@(43, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.114", "  +2.95%  ")
    ,@(44, "WEMIXToken", "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix", "2.31", "  +7.25%  ")
    ,@(45, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "2.46", "  +15.54%  ")
    ,@(46, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "2.093.78", "  +1.58%  ")
    ,@(47, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "3.33", "  +3.21%  ")
    ,@(48, "ApeXProtocol", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex", "2.23", "  +3.62%  ")
    ,@(49, "THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "5.56", "  +6.29%  ")
    ,@(50, "MultiversX", "https://coinranking.com/coin/omwkOTglq+multiversx-egld", "60.43", "  +2.10%  ")
    ,@(51, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "8.91", "  -0.45%  ")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Column D (Price) often looks numeric (e.g. "333.59"); Excel would
    # auto-convert such strings to a float. Force text storage, matching
    # the source workbook which stores these as plain text, then strip the
    # temporary text number-format so no stray cell style is left behind.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]
    $dCell.ClearFormats()

    $ws.Cells.Item($r, 5).Value = $row[4]
}

Write-Output "Done updating $($data.Count) rows"